# Auto-generated script updating F-column ("想去人数" / interest counts)
# values across all four worksheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4   # was 1
$ws.Range("F5").Value = 6163   # was 6154
$ws.Range("F6").Value = 685   # was 680
$ws.Range("F8").Value = 42   # was 37
$ws.Range("F9").Value = 98   # was 95
$ws.Range("F10").Value = 286   # was 283
$ws.Range("F12").Value = 621   # was 619
$ws.Range("F13").Value = 1073   # was 1070
$ws.Range("F16").Value = 324   # was 319
$ws.Range("F17").Value = 1390   # was 1387
$ws.Range("F18").Value = 622   # was 618
$ws.Range("F19").Value = 360   # was 358
$ws.Range("F20").Value = 81   # was 80
$ws.Range("F21").Value = 1040   # was 1037
$ws.Range("F22").Value = 91   # was 88
$ws.Range("F23").Value = 2106   # was 2105
$ws.Range("F25").Value = 59   # was 55
$ws.Range("F26").Value = 374   # was 373
$ws.Range("F28").Value = 3421   # was 3415

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 114   # was 113
$ws.Range("F9").Value = 674   # was 672
$ws.Range("F12").Value = 1012   # was 1011
$ws.Range("F21").Value = 4072   # was 4070
$ws.Range("F26").Value = 225   # was 224

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1535   # was 1534
$ws.Range("F10").Value = 124   # was 123
$ws.Range("F12").Value = 703   # was 702

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1535   # was 1534
$ws.Range("F8").Value = 124   # was 123
$ws.Range("F9").Value = 703   # was 702
$ws.Range("F10").Value = 114   # was 113
$ws.Range("F13").Value = 6163   # was 6154
$ws.Range("F15").Value = 685   # was 680
$ws.Range("F17").Value = 674   # was 672
$ws.Range("F18").Value = 42   # was 37
$ws.Range("F19").Value = 98   # was 95
$ws.Range("F20").Value = 286   # was 283
$ws.Range("F22").Value = 621   # was 619
$ws.Range("F27").Value = 1073   # was 1070
$ws.Range("F30").Value = 324   # was 319
$ws.Range("F33").Value = 1390   # was 1387
$ws.Range("F34").Value = 622   # was 618
$ws.Range("F35").Value = 360   # was 358
$ws.Range("F38").Value = 225   # was 224
$ws.Range("F39").Value = 1040   # was 1037
$ws.Range("F40").Value = 91   # was 88
$ws.Range("F42").Value = 2106   # was 2105
$ws.Range("F45").Value = 59   # was 55
$ws.Range("F46").Value = 374   # was 373
$ws.Range("F48").Value = 3421   # was 3415
